$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.541.46"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.62%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.877.45"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.06%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.04%  "

$ws.Range("E6").Value = "  -0.09%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4759"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.17%  "

$ws.Range("E8").Value = "  +1.46%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06506"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.15%  "

$ws.Range("E10").Value = "  +3.64%  "

$ws.Range("E11").Value = "  -0.29%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7395"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.70%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "96.63"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.53%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.874.33"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.11%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.194"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.36%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "274.54"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.35%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.628.67"
$ws.Range("D17").Style = "Normal"

$ws.Range("E18").Value = "  -1.20%  "

$ws.Range("E19").Value = "  -0.02%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007521"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.17%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.120.50"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.27%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9997"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.16%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.232"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.24%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.201"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.76%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "165.56"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.73%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.182"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.31%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.83"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.38%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.908"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.22%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.09848"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.10%  "

$ws.Range("E30").Value = "  -2.32%  "

$ws.Range("E31").Value = "  -0.39%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.269"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.76%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.104"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.64%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04814"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.38%  "

$ws.Range("E35").Value = "  +0.27%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6960"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.36%  "

$ws.Range("E37").Value = "  -0.05%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01869"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.77%  "

$ws.Range("E39").Value = "  +0.77%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.264"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.18%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "73.30"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.10%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.982"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.10%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4215"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.27%  "

$ws.Range("E45").Value = "  -0.73%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "101.92"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.10%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.387"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.29%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.992"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.22%  "

$ws.Range("E49").Value = "  +0.55%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "912.18"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.06%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05677"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.61%  "
